# Update "想去人数" (F4) and "想去人数" (F12) counts on the sheets that
# carry the 合肥 con listings: "展览" and "全部类型" mirror the same rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 4635
    $ws.Range("F12").Value = 188
}
